# Update "paises" (countries) workbook:
#  - refresh the "Datos actualizados" timestamp in A1
#  - update case counts for several countries whose row position is unchanged
#  - promote Albania / Uganda / Trinidad yTobago to a higher rank (their case
#    counts grew enough to overtake neighbouring countries in the
#    descending-by-total-cases ordering), pushing the countries that used to
#    occupy those ranks down by one row each

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 16:48"

# --- helper: write a full data row (B..H) ------------------------------
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- countries whose numbers update in place ---------------------------
# Estados Unidos (row 4)
Set-Row 4 5961926 6198 3256064 2523228 0 230 182634

# Argentina (row 15)
Set-Row 15 359638 0 268801 83176 0 98 7661

# Serbia (row 69)
Set-Row 69 30974 154 29343 924 0 2 707

# Republica de Macedonia (row 84)
Set-Row 84 13914 115 10507 2829 0 5 578

# Noruega (row 90)
Set-Row 90 10486 32 9150 1072 0 0 264

# Birmania (row 169)
Set-Row 169 580 76 341 233 0 0 6

# --- Albania overtakes Guayana Francesa & Croacia (rows 95-97) --------
$ws.Cells.Item(95, 1).Value = "Albania"
Set-Row 95 8927 168 4633 4031 0 4 263

$ws.Cells.Item(96, 1).Value = "Guayana Francesa"
Set-Row 96 8904 0 8399 449 0 0 56

$ws.Cells.Item(97, 1).Value = "Croacia"
Set-Row 97 8888 358 6362 2351 0 0 175

# --- Uganda overtakes Sudan del Sur (rows 132-133) ---------------------
$ws.Cells.Item(132, 1).Value = "Uganda"
Set-Row 132 2524 98 1268 1230 0 1 26

$ws.Cells.Item(133, 1).Value = "Sudan del Sur"
Set-Row 133 2507 0 1290 1170 0 0 47

# --- Trinidad yTobago overtakes Letonia, Burkina Faso, Togo, Liberia,
#     Reunion (rows 152-157) --------------------------------------------
$ws.Cells.Item(152, 1).Value = "Trinidad yTobago"
Set-Row 152 1384 132 178 1191 0 0 15

$ws.Cells.Item(153, 1).Value = "Letonia"
Set-Row 153 1360 18 1135 192 0 0 33

$ws.Cells.Item(154, 1).Value = "Burkina Faso"
Set-Row 154 1352 14 1058 239 0 0 55

$ws.Cells.Item(155, 1).Value = "Togo"
Set-Row 155 1309 0 919 363 0 0 27

$ws.Cells.Item(156, 1).Value = "Liberia"
Set-Row 156 1295 0 821 392 0 0 82

$ws.Cells.Item(157, 1).Value = "Reunion"
Set-Row 157 1292 0 692 594 0 0 6

Write-Host "edits applied"
